$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a few existing cells (Q57, O407, O410, R409, R410) ---
$ws.Cells.Item(57, 17).Value = 0      # Q57: 1 -> 0
$ws.Cells.Item(407, 15).Value = 2     # O407: 0 -> 2
$ws.Cells.Item(409, 18).Value = 0     # R409: (blank) -> 0
$ws.Cells.Item(410, 15).Value = 1     # O410: 0 -> 1
$ws.Cells.Item(410, 18).Value = 0     # R410: (blank) -> 0

# --- Append 5 new weekly rows (411-415) ---
$newRows = @(
    @{ Row=411; A=45474; B=265;               C=270.2999877929688; D=255;               E=262.75;             F=261.0379028320312; G=33356697; H=2024; I=7; J=1;  N=27; O=0; P=0; Q=0 },
    @{ Row=412; A=45481; B=262.6499938964844; C=263.75;             D=239.1000061035156; E=246.0500030517578; F=244.4467163085938; G=35211800; H=2024; I=7; J=8;  N=28; O=0; P=0; Q=0 },
    @{ Row=413; A=45488; B=246.4499969482422; C=248;                D=239.0500030517578; E=240.3500061035156; F=238.7838592529297; G=10852594; H=2024; I=7; J=15; N=29; O=0; P=1; Q=1 },
    @{ Row=414; A=45495; B=243.1999969482422; C=247;                D=227.1000061035156; E=235.8500061035156; F=234.3131866455078; G=58159770; H=2024; I=7; J=22; N=30; O=0; P=0; Q=0 },
    @{ Row=415; A=45502; B=238.3999938964844; C=241.6499938964844;  D=234.1000061035156; E=235.1999969482422; F=235.1999969482422; G=29536751; H=2024; I=7; J=29; N=31; O=0; P=0; Q=0 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A     # Datetime (serial date)
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $r.B     # Open
    $ws.Cells.Item($row, 3).Value = $r.C     # High
    $ws.Cells.Item($row, 4).Value = $r.D     # Low
    $ws.Cells.Item($row, 5).Value = $r.E     # Close
    $ws.Cells.Item($row, 6).Value = $r.F     # Adj Close
    $ws.Cells.Item($row, 7).Value = $r.G     # Volume
    $ws.Cells.Item($row, 8).Value = $r.H     # Year
    $ws.Cells.Item($row, 9).Value = $r.I     # Month
    $ws.Cells.Item($row, 10).Value = $r.J    # Day
    $ws.Cells.Item($row, 11).Value = 0       # Hour
    $ws.Cells.Item($row, 12).Value = 0       # Minute
    $ws.Cells.Item($row, 13).Value = 0       # Second
    $ws.Cells.Item($row, 14).Value = $r.N    # Week
    $ws.Cells.Item($row, 15).Value = $r.O    # isPivot
    $ws.Cells.Item($row, 16).Value = $r.P    # two_line_structure
    $ws.Cells.Item($row, 17).Value = $r.Q    # detect_structure
    # column R ("backup") is left blank for these new rows, matching the
    # source data (no backup value computed yet for un-reviewed weeks)
}
